# founda quad encoder and schmitt trigger inverter for debounce circuits
#
# Adds 5 new parts (rows 13-17) to the parts list, bumps the quantity of
# the existing 0-ohm resistor row (D10: 1 -> 9) and widens the
# Description column to fit the new longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# More 0-ohm resistors needed for the new debounce circuits.
$ws.Range("D10").Value = 9

# Seed rows 13-17 by cloning the formatting of the last existing data row
# (row 12) so the new rows pick up the same styles (Part Number wrap
# style + Hyperlink style) instead of minting new ones.
$ws.Range("A12:D12").Copy()
$ws.Range("A13:D17").PasteSpecial(-4122)

# Row 13 - 18pF cap for the 32.768kHz crystal
$ws.Range("A13").Value = "08055A180JAT2A"
$ws.Range("B13").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/avx-corporation/08055A180JAT2A/478-1307-1-ND/564339","Digikey - 478-1307-1-ND")'
$ws.Range("C13").Value = "18pF"
$ws.Range("D13").Value = 2

# Row 14 - 32.768kHz crystal
$ws.Range("A14").Value = "FX135A-327"
$ws.Range("B14").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/fox-electronics/FX135A-327/631-1002-1-ND/1024707","Digikey - 631-1002-1-ND")'
$ws.Range("C14").Value = "32.768kHz"
$ws.Range("D14").Value = 1

# Row 15 - 20pF cap for the 32.768kHz crystal
$ws.Range("A15").Value = "08051A200JAT2A"
$ws.Range("B15").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/avx-corporation/08051A200JAT2A/478-3735-1-ND/1116433","Digikey - 478-3735-1-ND")'
$ws.Range("C15").Value = "20pF"
$ws.Range("D15").Value = 2

# Row 16 - quad encoder
$ws.Range("A16").Value = "EN11-HSM1BF20"
$ws.Range("B16").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/tt-electronics-bi/EN11-HSM1BF20/987-1398-ND/2620667","Digikey - 987-1398-ND")'
$ws.Range("C16").Value = "Quad Encoder"
$ws.Range("D16").Value = 1

# Row 17 - 6x schmitt trigger inverter, used for debounce
$ws.Range("A17").Value = "74HC14D,653"
$ws.Range("B17").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/nxp-semiconductors/74HC14D,653/568-1401-1-ND/763376","Digikey - 568-1401-1-ND")'
$ws.Range("C17").Value = "6x Schmitt trigger inverter"
$ws.Range("D17").Value = 1

# Widen the Description column to fit the new text.
$ws.Columns.Item(3).ColumnWidth = 22.3

# Leave the cursor where the author left it after adding the new rows.
$ws.Range("B16").Select() | Out-Null
